# Scheduled runner update: refresh computed profit/price figures across the
# per-job Leve tables (ALC, ARM, BSM, CRP, CUL, LTW, WVR). Only numeric value
# cells are touched; labels/formulas/styles are left untouched.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 786.08
$ws.Range("I15").Value = 786.08
$ws.Range("K15").Value = 2358.24
$ws.Range("M15").Value = -2189.24
$ws.Range("H69").Value = 4216.6665
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 4600
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 13800
$ws.Range("M69").Value = -6026
$ws.Range("N69").Value = -15548
$ws.Range("H72").Value = 4216.6665
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 4600
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 41400
$ws.Range("M72").Value = -16332
$ws.Range("N72").Value = -50136
$ws.Range("H112").Value = 3087462.2
$ws.Range("J112").Value = 3087462.2
$ws.Range("L112").Value = 9262386.600000001
$ws.Range("N112").Value = -9264602.600000001
$ws.Range("H129").Value = 218523.34
$ws.Range("J129").Value = 245132.05
$ws.Range("L129").Value = 735396.1499999999
$ws.Range("N129").Value = -745396.1499999999
$ws.Range("H133").Value = 49253
$ws.Range("J133").Value = 49253
$ws.Range("L133").Value = 49253
$ws.Range("N133").Value = -59373
$ws.Range("H135").Value = 10641168
$ws.Range("I135").Value = 429.9091
$ws.Range("J135").Value = 166705330
$ws.Range("K135").Value = 3869.1819
$ws.Range("L135").Value = 1500347970
$ws.Range("M135").Value = -1334.1819
$ws.Range("N135").Value = -1500353040
$ws.Range("H137").Value = 3595.9092
$ws.Range("I137").Value = 3944.125
$ws.Range("J137").Value = 2667.3333
$ws.Range("K137").Value = 11832.375
$ws.Range("L137").Value = 8001.999899999999
$ws.Range("M137").Value = -9282.375
$ws.Range("N137").Value = -13101.9999
$ws.Range("H138").Value = 15387827
$ws.Range("J138").Value = 3574.796
$ws.Range("L138").Value = 10724.388
$ws.Range("N138").Value = -21004.388
$ws.Range("H140").Value = 30780
$ws.Range("J140").Value = 30780
$ws.Range("L140").Value = 30780
$ws.Range("N140").Value = -41140

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4627.0513
$ws.Range("I32").Value = 3107.1287
$ws.Range("K32").Value = 3107.1287
$ws.Range("M32").Value = -2820.1287
$ws.Range("H61").Value = 368614.22
$ws.Range("I61").Value = 410086.3
$ws.Range("J61").Value = 3659.8
$ws.Range("K61").Value = 410086.3
$ws.Range("L61").Value = 3659.8
$ws.Range("M61").Value = -409874.3
$ws.Range("N61").Value = -4083.8
$ws.Range("H132").Value = 17702.719
$ws.Range("J132").Value = 51761.5
$ws.Range("L132").Value = 155284.5
$ws.Range("N132").Value = -160344.5
$ws.Range("H136").Value = 368614.22
$ws.Range("I136").Value = 410086.3
$ws.Range("J136").Value = 3659.8
$ws.Range("K136").Value = 1230258.9
$ws.Range("L136").Value = 10979.4
$ws.Range("M136").Value = -1227708.9
$ws.Range("N136").Value = -16079.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3183.5789
$ws.Range("I134").Value = 3438.6667
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 10316.0001
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -7781.000100000001
$ws.Range("N134").Value = -9570

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7965.6113
$ws.Range("I31").Value = 9627.2
$ws.Range("J31").Value = 7326.5386
$ws.Range("K31").Value = 9627.2
$ws.Range("L31").Value = 7326.5386
$ws.Range("M31").Value = -9332.2
$ws.Range("N31").Value = -7916.5386
$ws.Range("H34").Value = 7965.6113
$ws.Range("I34").Value = 9627.2
$ws.Range("J34").Value = 7326.5386
$ws.Range("K34").Value = 9627.2
$ws.Range("L34").Value = 7326.5386
$ws.Range("M34").Value = -9425.2
$ws.Range("N34").Value = -7730.5386
$ws.Range("H107").Value = 909.53125
$ws.Range("I107").Value = 515
$ws.Range("J107").Value = 1116.1904
$ws.Range("K107").Value = 515
$ws.Range("L107").Value = 1116.1904
$ws.Range("M107").Value = 1405
$ws.Range("N107").Value = -4956.190399999999
$ws.Range("H132").Value = 2143.7173
$ws.Range("I132").Value = 1658.2222
$ws.Range("J132").Value = 3891.5
$ws.Range("K132").Value = 4974.6666
$ws.Range("L132").Value = 11674.5
$ws.Range("M132").Value = -2444.6666
$ws.Range("N132").Value = -16734.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 615
$ws.Range("H131").Value = 776.09
$ws.Range("J131").Value = 795.84045
$ws.Range("L131").Value = 2387.52135
$ws.Range("N131").Value = -12467.52135

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4262.033
$ws.Range("I7").Value = 2759.75
$ws.Range("J7").Value = 5978.9287
$ws.Range("K7").Value = 2759.75
$ws.Range("L7").Value = 5978.9287
$ws.Range("M7").Value = -2647.75
$ws.Range("N7").Value = -6202.9287
$ws.Range("H46").Value = 1361.875
$ws.Range("I46").Value = 899.1667
$ws.Range("J46").Value = 2750
$ws.Range("K46").Value = 899.1667
$ws.Range("L46").Value = 2750
$ws.Range("M46").Value = -711.1667
$ws.Range("N46").Value = -3126
$ws.Range("H122").Value = 819845.9
$ws.Range("I122").Value = 1785182
$ws.Range("J122").Value = 3023
$ws.Range("K122").Value = 5355546
$ws.Range("L122").Value = 9069
$ws.Range("M122").Value = -5353096
$ws.Range("N122").Value = -13969
$ws.Range("H126").Value = 4262.033
$ws.Range("I126").Value = 2759.75
$ws.Range("J126").Value = 5978.9287
$ws.Range("K126").Value = 8279.25
$ws.Range("L126").Value = 17936.7861
$ws.Range("M126").Value = -5809.25
$ws.Range("N126").Value = -22876.7861
$ws.Range("H132").Value = 1618.1842
$ws.Range("I132").Value = 1449.7667
$ws.Range("J132").Value = 2249.75
$ws.Range("K132").Value = 4349.300099999999
$ws.Range("L132").Value = 6749.25
$ws.Range("M132").Value = -1819.300099999999
$ws.Range("N132").Value = -11809.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1105
$ws.Range("I126").Value = 1082.3077
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 3246.9231
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -776.9231
$ws.Range("N126").Value = -9140
$ws.Range("H132").Value = 726.7761
$ws.Range("I132").Value = 521.4259
$ws.Range("J132").Value = 1579.7693
$ws.Range("K132").Value = 1564.2777
$ws.Range("L132").Value = 4739.3079
$ws.Range("M132").Value = 965.7223000000001
$ws.Range("N132").Value = -9799.3079
$ws.Range("H136").Value = 15153376
$ws.Range("I136").Value = 24391192
$ws.Range("J136").Value = 3359.8
$ws.Range("K136").Value = 73173576
$ws.Range("L136").Value = 10079.4
$ws.Range("M136").Value = -73171026
$ws.Range("N136").Value = -15179.4
